$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values of the columns that get reshuffled (D, J, K, L, M, P)
# for each data row (2..15) before making any changes.
$cols = @("D", "J", "K", "L", "M", "P")
$orig = @{}
foreach ($r in 2..15) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $row
}

# Mapping of target row -> source row describing where each row's
# D/J/K/L/M/P values come from after the edit.
$mapping = @{
    2  = 15
    3  = 8
    4  = 11
    5  = 9
    6  = 13
    7  = 7
    8  = 10
    9  = 6
    10 = 12
    11 = 3
    12 = 5
    13 = 4
    14 = 14
    15 = 2
}

foreach ($targetRow in ($mapping.Keys | Sort-Object)) {
    $sourceRow = $mapping[$targetRow]
    if ($sourceRow -eq $targetRow) {
        continue
    }
    $src = $orig[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $src[$c]
    }
}
